$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 274; this shifts rows 274:335 down to 275:336
$ws.Rows.Item(274).Insert()

# Populate the newly inserted row 274 with its data (matching the
# surrounding rows' layout/template, since it is a new daily price record)
$ws.Range("A274").Value = 3
$ws.Range("B274").Value = "Femacal de La Calera"
$ws.Range("C274").Value = "Coquimbo"
$ws.Range("D274").Value = 44785
$ws.Range("E274").Value = 5
$ws.Range("F274").Value = 100112001
$ws.Range("G274").Value = "Berenjena"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 105
$ws.Range("K274").Value = 8500
$ws.Range("L274").Value = 9000
$ws.Range("M274").Value = 8738
$ws.Range("N274").Value = "$/caja 60 unidades"
$ws.Range("O274").Value = "Región de Arica y Parinacota"
$ws.Range("P274").Value = 146
$ws.Range("Q274").Value = 60
$ws.Range("R274").Value = "Hortaliza"
